# "reading metadata for library"
#
# The sheet originally held subcategory rows (A:B) plus some loosely-named
# "default category/subcategory" bookkeeping columns (C:H). This rewrites it
# to a clean subcategory -> category metadata table: every subcategory row
# now carries its owning category's key/name/sort-priority (read off the
# rows that already had that metadata), and the leftover "is default ..."
# flag columns are dropped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "Is default subcategory" column (old D) and the "Is default
# category" column (old H, which is G after the first delete). This also
# slides Category/Category name/Category priority left into D/E/F and
# shrinks the used range from H back to F.
$ws.Columns("D").Delete()
$ws.Columns("G").Delete()

# Remove the now-orphaned header label for column A (subcategory key) --
# the new header row only labels B:F.
$ws.Range("A1").ClearContents()

# Re-label the remaining headers with the new snake_case metadata names.
$ws.Range("B1").Value = "sheet_name"
$ws.Range("C1").Value = "category"
$ws.Range("D1").Value = "category_name"
$ws.Range("E1").Value = "subcategory_sort_priority"
$ws.Range("F1").Value = "category_sort_priority"

# The old sort-priority column (now E) keeps its per-row numbers as-is; only
# C/D/F need the category metadata filled in (rows 3-6 previously had it
# blank). Do this per-row so every subcategory ends up fully tagged with its
# category's key, display name, and sort priority.
function Set-CategoryMeta($row, $category, $categoryName, $categoryPriority) {
    $ws.Cells.Item($row, 3).Value = $category
    $ws.Cells.Item($row, 4).Value = $categoryName
    $ws.Cells.Item($row, 6).Value = $categoryPriority
}

Set-CategoryMeta 2 "Прочее" "Прочие материалы" 3
Set-CategoryMeta 3 "Прочее" "Прочие материалы" 3
Set-CategoryMeta 4 "Прочее" "Прочие материалы" 3
Set-CategoryMeta 5 "Прочее" "Прочие материалы" 3
Set-CategoryMeta 6 "Прочее" "Прочие материалы" 3
Set-CategoryMeta 7 "Кабель" "Кабельные изделия" 2
Set-CategoryMeta 8 "Свет" "Светотехническое оборудование" 1

# subcategory_sort_priority (E) values are unchanged from the original
# sort-priority column, just re-set explicitly so the column is unambiguous.
$ws.Cells.Item(2, 5).Value = 5
$ws.Cells.Item(3, 5).Value = 4
$ws.Cells.Item(4, 5).Value = 4
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(8, 5).Value = 0

Write-Output "done"
